# ARLIS_contents.xlsx edit script
# Applies: new "Added just to get examples out" column (D) marking all existing
# predicate/argument pairs as TRUE, updates the bribe/GAIN sentiment row, and
# appends new lexicon rows for sovereignty / slavery / project, replacing the
# old scratch comment string.
#
# NOTE: the order in which brand-new text values are first written matters,
# since the saved shared-string table lists new strings in first-write order
# (appended after the strings that already existed). To reproduce the target
# table order we write, in order: the new D1 header text, then "sovereignty",
# then "slavery", then "project".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D header: "Added just to get examples out" (write this first) ---
$ws.Range("A2").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Added just to get examples out"

# --- Existing rows 2-10: text / numeric values are unchanged by the diff      ---
# --- (only shared-string indices shifted because of insertions/removals), so ---
# --- we leave A/B/C alone except for row 8 (bribe) whose sentiment flips.    ---

# Row 8 ("bribe"/GAIN): sentiment flips from -1.0 to 1.0, and its style moves
# from the "normal" data style to the "special" style already used by D8.
$ws.Range("D8").Copy()
$ws.Range("C8").PasteSpecial(-4122)   # xlPasteFormats: adopt D8's style (s=3)
$ws.Range("C8").Value = 1.0

# Add TRUE marker in column D for every existing predicate/argument row,
# reusing the format already present on neighboring column-A cells (s=2).
$rows2to10 = 2,3,4,5,6,7,9,10
foreach ($r in $rows2to10) {
    $srcCell = "A$r"
    $dstCell = "D$r"
    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial(-4122)
    $ws.Range($dstCell).Value = $true
}
# D8 already exists with style s=3; just change its value from the old
# scratch comment to the boolean TRUE marker, keeping its style.
$ws.Range("D8").Value = $true

# --- New rows 11-15 ---

# Row 11: sovereignty / LOSE / -1.0 / TRUE (D11 uses the "special" style)
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "sovereignty"

$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "LOSE"

$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = -1.0

$ws.Range("D8").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = $true

# Row 12: sovereignty / GAIN / 1.0 (special style) / TRUE
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "sovereignty"

$ws.Range("B3").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = "GAIN"

$ws.Range("D8").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").Value = 1.0

$ws.Range("A2").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D12").Value = $true

# Row 13: slavery / LOSE / -1.0 / TRUE
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "slavery"

$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "LOSE"

$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").Value = -1.0

$ws.Range("A2").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D13").Value = $true

# Row 14: project / LOSE / -1.0 / TRUE
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "project"

$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "LOSE"

$ws.Range("C2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = -1.0

$ws.Range("A2").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = $true

# Row 15: project / GAIN (special style) / 1.0 (special style) / TRUE
$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "project"

$ws.Range("D8").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "GAIN"

$ws.Range("D8").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1.0

$ws.Range("A2").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = $true
